# Add new P2P member "Bethany Reinhart" (Canada / Fisheries and Ocean Canada)
# as a new row 9, pushing all subsequent members down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 9 ("Brian Helmuth"), shifting
# everything from row 9 down onward to row 10 onward.
$ws.Rows("9:9").Insert()

# Populate the new row with the new member's data.
$ws.Range("A9").Value = "Bethany Reinhart"
$ws.Range("B9").Value = "Canada"
$ws.Range("C9").Value = "[Fisheries and Ocean Canada](https://www.mar.dfo-mpo.gc.ca/SABS/Home)"
$ws.Range("D9").Value = "RS/SB"
$ws.Range("E9").Value = "bethany.reinhart@dfo-mpo.gc.ca"

# Match the saved selection/view state of the edited workbook.
$ws.Range("E10").Select() | Out-Null
